# Weekly update: add a new Perejil (Vega Monumental Concepción) price entry
# dated 2021-10-19 (serial 44488) at the top of the data block (rows 60-61),
# pushing all existing rows down by two. The new pair duplicates the
# (previously topmost) row 60/61 record in every column except the date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing row 60:61 pair (Primera / Segunda) before shifting -
# the new rows will reuse this data, only the date changes.
$srcRange = $ws.Range("A60:R61")
$srcValues = $srcRange.Value()

# Insert two new blank rows at row 60, shifting rows 60:109 down to 62:111.
$ws.Range("A60:A61").EntireRow.Insert()

# Re-populate the newly inserted rows 60:61 with the captured data.
$ws.Range("A60:R61").Value = $srcValues

# New week's date (2021-10-19) for the newly inserted entries.
$ws.Range("D60").Value = 44488
$ws.Range("D61").Value = 44488
